$wb = $excel.ActiveWorkbook

# Sheet "1-15": just update the Month/Year label, no layout change
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A1").Value = "Month/Year: FEBRUARY 2021"

# Sheet "16-End": update the label, and shrink the day columns from
# O (31) down to M (28) since February only has 13 days left after the 16th
# of the month (16..28), instead of 15 (16..30).
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "Month/Year: FEBRUARY 2021"

# Carry the closing border formatting from the old last column (O1) onto
# the new last column (M1) before removing the now-unused N:O columns.
$ws2.Range("O1").Copy()
$ws2.Range("M1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the now unused N:O columns; Excel will shrink the merged header
# range and the sheet dimension automatically.
$ws2.Range("N1:O1").EntireColumn.Delete()
